# Generate Report for Handoff
# Regenerates the localization-status report for the new source file
# (cedffbac-5c2e-4c4d-b936-9510fcd941bc.md, which replaces
# 05e3eb0f-55c7-48b4-8b4f-120ac4583c02.md) and refreshes the handoff
# timestamps / handoff xliff artifacts on the Overview, zh-cn and de-de
# sheets.

$wb = $excel.ActiveWorkbook

$oldName = "05e3eb0f-55c7-48b4-8b4f-120ac4583c02"
$newName = "cedffbac-5c2e-4c4d-b936-9510fcd941bc"
$newHash = "6b960fb1a4d4edcf92332fcee6302c084842363e"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newName.md"

$wsOverview.Range("B2").Value = "e2e\$newName.md"
foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newName.md"
    }
}

$wsOverview.Range("G2").Value = "2016-08-25 02:59:22"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newName.md"
foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newName.md"
    }
}

$wsZhCn.Range("G2").Value = "$newName.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-25 02:59:17"

# Latest Target File (I2) loses its handback hyperlink/value entirely.
$toDeleteZhCn = @()
foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $toDeleteZhCn += $hl
    }
}
foreach ($hl in $toDeleteZhCn) {
    $hl.Delete()
}
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").ClearFormats()

# Latest Handback File (J2) is cleared too.
$wsZhCn.Range("J2").Value = ""

# Latest Handback DateTime (K2) resets to the zero-date sentinel.
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Column widths for I/J shrink from the generic 40 placeholder to their
# content-fit sizes.
$wsZhCn.Columns.Item(9).ColumnWidth = 17.8
$wsZhCn.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newName.md"
foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newName.md"
    }
}

$wsDeDe.Range("G2").Value = "$newName.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-25 02:59:22"

# Latest Target File (I2) loses its handback hyperlink/value entirely.
$toDeleteDeDe = @()
foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $toDeleteDeDe += $hl
    }
}
foreach ($hl in $toDeleteDeDe) {
    $hl.Delete()
}
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").ClearFormats()

# Latest Handback File (J2) is cleared too.
$wsDeDe.Range("J2").Value = ""

# Latest Handback DateTime (K2) resets to the zero-date sentinel.
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Columns.Item(9).ColumnWidth = 17.8
$wsDeDe.Columns.Item(10).ColumnWidth = 20.8
